$d = $word.ActiveDocument

# Map of old numeric suffix -> new numeric suffix for each inline picture's
# alt-text / description (drives both wp:docPr@descr and pic:cNvPr@descr,
# since Word's InlineShape.AlternativeText covers both at once).
$suffixMap = @(
    @{ Old = "85330"; New = "10463" },
    @{ Old = "85335"; New = "10468" },
    @{ Old = "85348"; New = "104621" },
    @{ Old = "85352"; New = "104625" },
    @{ Old = "85356"; New = "104629" }
)

for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $oldText = $shape.AlternativeText
    $mapping = $suffixMap[$i - 1]
    $newText = $oldText -replace ($mapping.Old + "\.jpg$"), ($mapping.New + ".jpg")
    $shape.AlternativeText = $newText
}
